$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the June 16th (2020-06-16, serial 43998) SSA raw/clean row ---

# The new date cell (B17) should get the "date-only" format (style 3) that
# B16 currently has, so grab that format first...
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)  # xlPasteFormats

# ...then B16 itself switches to the "date + time" format (style 2) used by
# all the other date cells in the column (B2:B15).
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)  # xlPasteFormats

# A17 reuses the bordered/bold index-column format already applied to A16.
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Row 17 values: raw (C:F) and clean (G) SSA figures for 2020-06-16.
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 43998
$ws.Range("C17").Value = 154863
$ws.Range("D17").Value = 216857
$ws.Range("E17").Value = 56843
$ws.Range("F17").Value = 18310
$ws.Range("G17").Value = 32.29
